# Generación de oficio de Investigación con Policía Ministerial.
# Applies the field-level corrections to FormatoDenuncia1.docx.

$d = $word.ActiveDocument

# --- Encabezado / carpeta -------------------------------------------------
# Tabla 1: FECHA DE INICIO
$d.Tables(1).Cell(3, 2).Range.Text = "13/12/2017"

# --- Datos del denunciante (Tabla 2) ---------------------------------------
$d.Tables(2).Cell(1, 2).Range.Text = "CASA AMERICA  "
$d.Tables(2).Cell(4, 2).Range.Text = "ALEXANDER VON HUMBOLT #2, COLONIA ACAJETE, ACAJETE, VERACRUZ"
$d.Tables(2).Cell(5, 2).Range.Text = "SI"
$d.Tables(2).Cell(6, 2).Range.Text = "CASAAMERICA001"
$d.Tables(2).Cell(12, 2).Range.Text = "ALEXANDER VON HUMBOLT #2, COLONIA ACAJETE, ACAJETE, VERACRUZ"

# --- Domicilio / contacto del denunciante (Tabla 3) -------------------------
$d.Tables(3).Cell(1, 2).Range.Text = "ALEXANDER VON HUMBOLt #2, COLONIA ACAJETE, ACAJETE, VERACRUZ"
$d.Tables(3).Cell(2, 2).Range.Text = "sebastianlobato20@gmail.com"
$d.Tables(3).Cell(3, 2).Range.Text = "123456789"
$d.Tables(3).Cell(3, 4).Range.Text = "123456789"

# --- Lugar de los hechos (Tabla 4) ------------------------------------------
$d.Tables(4).Cell(1, 2).Range.Text = "ALEXANDER VON HUMBOLT #1, COLONIA ACAJETE, ACAJETE, VERACRUZ"
$d.Tables(4).Cell(2, 2).Range.Text = "CRISTAL"
$d.Tables(4).Cell(2, 6).Range.Text = "ACAJETE"
$d.Tables(4).Cell(3, 2).Range.Text = "13/12/2017 21:59:00"
$d.Tables(4).Cell(4, 2).Range.Text = "ASDFGHJKLÑ"
$d.Tables(4).Cell(4, 5).Range.Text = "QWERTYUIOP"

# --- Delito (Tabla 5) --------------------------------------------------------
$d.Tables(5).Cell(1, 2).Range.Text = "ROBO A NEGOCIACIONES SIN VIOLENCIA"
$d.Tables(5).Cell(1, 4).Range.Text = "NO"
$d.Tables(5).Cell(3, 2).Range.Text = "CULPOSO"

# --- Presunto responsable (Tabla 6) ------------------------------------------
$d.Tables(6).Cell(1, 2).Range.Text = "QUIEN RESULTE RESPONSABLE  "
$d.Tables(6).Cell(2, 2).Range.Text = "SIN INFORMACION #S/N, COLONIA SIN INFORMACION, SIN INFORMACION, SIN INFORMACION"
$d.Tables(6).Cell(3, 2).Range.Text = "NO"
$d.Tables(6).Cell(3, 4).Range.Text = "SIN INFORMACION"

# --- Narrativa de los hechos y firma (Tabla 8) -------------------------------
$d.Content.Find.Execute("Se llevaron mi nueva versión de android", $true, $false, $false, $false, $false, $true, 1, $false, "ENTRARON A MI TIENDA A LAS 2 DE LA MADRUGADA Y SE LLEVARON TODITO.", 2)
$d.Content.Find.Execute("LIC. NAIN LOBATO GARCíA", $true, $false, $false, $false, $false, $true, 1, $false, "LIC. NAIN LOBATO GARCÍA", 2)
